$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values would otherwise be auto-interpreted as
# numbers by Excel (e.g. "312.39") must be forced to a Text number
# format first, so they are written back out as literal strings,
# matching the original inline-string cell type/content exactly.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.061.04"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.830.08"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").Value = "312.39"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "0.4601"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").Value = "0.3701"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").Value = "0.07344"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "0.8717"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").Value = "0.07928"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "19.81"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "1.826.89"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "5.345"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "6.548"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "91.79"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "0.000008875"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "26.920.17"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("D22").Value = "5.123"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "2.004.78"
$ws.Range("E24").Value = "  -4.42%  "
$ws.Range("D25").Value = "152.54"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").Value = "2.073"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "5.126"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "115.19"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("D31").Value = "0.08870"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "2.975"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "0.7324"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "4.438"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").Value = "2.452"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").Value = "1.073"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "0.05241"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "0.01936"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "2.949"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "7.161"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").Value = "0.5155"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "0.8825"
$ws.Range("E43").Value = "  -12.88%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.1631"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "8.246"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "0.4833"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "10.25"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "102.30"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "1.627"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "0.06218"
$ws.Range("E51").Value = "  -0.99%  "
